# guide41_email.xlsx — "Add files via upload" edit
#
# Net effect (derived from the OOXML diff):
#   * Sheet "index" content is unchanged (only shared-string indices shift
#     because strings are removed elsewhere) — nothing to do there except
#     move the active selection / active sheet.
#   * Sheet "p1" drops its leading "description" placeholder row, so every
#     later row shifts up by one; its header link B1 becomes a plain title
#     string ("登録メールアドレスを変更する") instead of the old
#     phonetic "登録メールアドレス"; and the big HTML paragraph row +
#     the "<h3>...email change</h3>" anchor row are removed and replaced
#     by the mail2.png screenshot row that used to be further down.
#   * "p1" becomes the active tab, scrolled/selected to B6, zoomed to 85%.

$wb = $excel.ActiveWorkbook
$wsIndex = $wb.Worksheets.Item("index")
$wsP1 = $wb.Worksheets.Item("p1")

# --- p1: remove the now-unused "description" header row ------------------
# (A1:"description", B1: blank) — everything below shifts up by one row.
$wsP1.Rows(1).Delete()

# --- p1: remove the old long description paragraph + "email change"
#     anchor rows (now rows 5 and 6 after the shift above) -- the mail1/
#     mail2/mail3 screenshot rows that followed shift up to take their
#     place (old rows 8/9/10 -> new rows 5/6/7).
$wsP1.Rows("5:6").Delete()

# --- p1: the page header text changes from the old phonetic
#     "登録メールアドレス" label to the new plain title string.
$wsP1.Range("B1").Value = "登録メールアドレスを変更する"

# --- view state: "p1" becomes the active/selected sheet, scrolled to B6
#     at 85% zoom; "index" keeps its own cursor position at B6 too.
$wsIndex.Activate()
$wsIndex.Range("B6").Select()

$wsP1.Activate()
$excel.ActiveWindow.Zoom = 85
$wsP1.Range("B6").Select()
